$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Rename some customers to their acronyms
$ws.Range("A22").Value = "DDL"
$ws.Range("A4").Value = "CMML"
$ws.Range("A7").Value = "KSFH"
$ws.Range("A6").Value = "SHCH"

# Update coordinates for row 24 (more precise values)
$ws.Range("B24").Value = 104.132008743602
$ws.Range("C24").Value = 10.570141163026801

# Add a new customer row
$ws.Range("A31").Value = "Gold Medical Diagnostic Laboratory"
$ws.Range("B31").Value = 104.926276
$ws.Range("C31").Value = 11.572371
$ws.Range("C31").Style = "Normal"

# Update the view: scroll down a bit and move the selection
$ws.Activate()
$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B35").Select()
